$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C ("Type"), shifting SiO2_liq..T right by one.
$ws.Range("C1").EntireColumn.Insert()

# New header
$ws.Range("C1").Value = "Type"

# New values for data rows 2-6
$ws.Range("C2:C6").Value = "Matrix"

# Update selection to match target state
$ws.Range("D13").Select()
